# Daily update at 8 AM UTC
# Adds the next day's row (row 62) to the Wins_Over_Time tracking sheet,
# and moves the "last row" date number format from the old last row (61)
# to the new last row (62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the special "last row" number format (currently on A61) before
# we overwrite anything, then give A61 the plain date format every other
# data row uses (copied from A60, the row above it).
$ws.Range("A62").NumberFormat = $ws.Range("A61").NumberFormat
$ws.Range("A61").NumberFormat = $ws.Range("A60").NumberFormat

# Append the new day's data in row 62.
$ws.Range("A62").Value = 45648
$ws.Range("B62").Value = 144
$ws.Range("C62").Value = 134
$ws.Range("D62").Value = 142
